$wb = $excel.ActiveWorkbook

# The "Status" value shrank from "Ready for handoff" to "In Translation",
# which is referenced from the Overview summary sheet (zh-cn/de-de columns)
# as well as the per-locale status sheets. Updating the shared text causes
# Excel to re-autofit the (now narrower) status columns.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C:C").ColumnWidth = 12.5
